$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)

$shp = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$shp.Name = "Shape 89"

$shp.Left = 470.7874315748031
$shp.Top = 359.13782527559056
$shp.Width = 205.15743157480316
$shp.Height = 26.10239220472441

$shp.Fill.Visible = $false
$shp.Line.Visible = $false

$tf = $shp.TextFrame
$tf.AutoSize = 0
$tf.MarginLeft = 7.198848897637795
$tf.MarginRight = 7.198848897637795
$tf.MarginTop = 7.198848897637795
$tf.MarginBottom = 7.198848897637795
$tf.VerticalAnchor = 1

$tr = $tf.TextRange
$tr.Text = "(http://littlegreenriver.com/weblog/wp-content/uploads/mtv-diagram-730x1024.png)"
$tr.Font.Size = 8
$tr.Font.Color.RGB = 13421772

$pf = $tr.ParagraphFormat
$pf.SpaceBefore = 0
$pf.Bullet.Visible = $false

Write-Output "done"
